$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.910.90"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "3.360.33"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'565.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "'148.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'7.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "3.936.90"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "'27.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "3.366.16"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").Value = "61.002.19"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").Value = "'8.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("D21").Value = "'374.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").Value = "'75.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").Value = "'0.558"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "3.501.34"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "'0.0000108"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.62%  "
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("D33").Value = "'22.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("E34").Value = "  -4.80%  "
$ws.Range("D35").Value = "'5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'169.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "'29.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.04%  "
$ws.Range("D40").Value = "3.396.47"
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").Value = "'1.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("D47").Value = "2.488.37"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").Value = "'6.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'22.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("E51").Value = "  -2.71%  "
